$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1691.7368
$ws.Range("I29").Value = 703.44446
$ws.Range("J29").Value = 2581.2
$ws.Range("K29").Value = 2110.33338
$ws.Range("L29").Value = 7743.599999999999
$ws.Range("M29").Value = -1829.33338
$ws.Range("N29").Value = -8305.599999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 12666
$ws.Range("I43").Value = 12666
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 12666
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -12597

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 750
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 750
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2250
$ws.Range("N80").Value = -4246

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 750
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 750
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 6750
$ws.Range("N83").Value = -16734

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 501.375
$ws.Range("I92").Value = 443.66666
$ws.Range("J92").Value = 674.5
$ws.Range("K92").Value = 443.66666
$ws.Range("L92").Value = 674.5
$ws.Range("M92").Value = 804.33334
$ws.Range("N92").Value = -3170.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2985.4167
$ws.Range("I132").Value = 1438.6364
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 4315.9092
$ws.Range("L132").Value = 60000
$ws.Range("M132").Value = -1785.9092
$ws.Range("N132").Value = -65060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6645.3335
$ws.Range("I61").Value = 6645.3335
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6645.3335
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6433.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1370.0605
$ws.Range("I110").Value = 936.7778
$ws.Range("J110").Value = 1890
$ws.Range("K110").Value = 936.7778
$ws.Range("L110").Value = 1890
$ws.Range("M110").Value = 1108.2222
$ws.Range("N110").Value = -5980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2580.652
$ws.Range("I122").Value = 2009.6875
$ws.Range("J122").Value = 3885.7144
$ws.Range("K122").Value = 6029.0625
$ws.Range("L122").Value = 11657.1432
$ws.Range("M122").Value = -3579.0625
$ws.Range("N122").Value = -16557.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2738.6667
$ws.Range("I132").Value = 2644
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 7932
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -5402
$ws.Range("N132").Value = -14057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6645.3335
$ws.Range("I136").Value = 6645.3335
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 19936.0005
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -17386.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 55000
$ws.Range("I69").Value = 55000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 55000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -54189

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 55000
$ws.Range("I72").Value = 55000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 165000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -160944

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 626322.0600000001
$ws.Range("I86").Value = 1236.2727
$ws.Range("J86").Value = 2001510.8
$ws.Range("K86").Value = 1236.2727
$ws.Range("L86").Value = 2001510.8
$ws.Range("M86").Value = -113.2727
$ws.Range("N86").Value = -2003756.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 626322.0600000001
$ws.Range("I89").Value = 1236.2727
$ws.Range("J89").Value = 2001510.8
$ws.Range("K89").Value = 6181.363499999999
$ws.Range("L89").Value = 10007554
$ws.Range("M89").Value = -565.3634999999995
$ws.Range("N89").Value = -10018786

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1557.6842
$ws.Range("I134").Value = 1557.6842
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4673.0526
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2138.0526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 20499
$ws.Range("I3").Value = 19249.5
$ws.Range("J3").Value = 22998
$ws.Range("K3").Value = 19249.5
$ws.Range("L3").Value = 22998
$ws.Range("M3").Value = -19136.5
$ws.Range("N3").Value = -23224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1480.3572
$ws.Range("I7").Value = 1556.6364
$ws.Range("J7").Value = 1200.6666
$ws.Range("K7").Value = 1556.6364
$ws.Range("L7").Value = 1200.6666
$ws.Range("M7").Value = -1443.6364
$ws.Range("N7").Value = -1426.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 129840.5
$ws.Range("I22").Value = 255499.75
$ws.Range("J22").Value = 4181.25
$ws.Range("K22").Value = 255499.75
$ws.Range("L22").Value = 4181.25
$ws.Range("M22").Value = -255149.75
$ws.Range("N22").Value = -4881.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 69999
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 69999
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 69999
$ws.Range("N68").Value = -71497

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 69999
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 69999
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 209997
$ws.Range("N71").Value = -217485

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 628.5625
$ws.Range("I105").Value = 658.3077
$ws.Range("J105").Value = 499.66666
$ws.Range("K105").Value = 658.3077
$ws.Range("L105").Value = 499.66666
$ws.Range("M105").Value = 1088.6923
$ws.Range("N105").Value = -3993.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 63749.75
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 63749.75
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 63749.75
$ws.Range("N109").Value = -65829.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2061.8823
$ws.Range("I134").Value = 1561.9231
$ws.Range("J134").Value = 3686.75
$ws.Range("K134").Value = 4685.7693
$ws.Range("L134").Value = 11060.25
$ws.Range("M134").Value = -2150.7693
$ws.Range("N134").Value = -16130.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 92584
$ws.Range("I23").Value = 1743.6666
$ws.Range("J23").Value = 126649.125
$ws.Range("K23").Value = 5230.9998
$ws.Range("L23").Value = 379947.375
$ws.Range("M23").Value = -4995.9998
$ws.Range("N23").Value = -380417.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 8665.333000000001
$ws.Range("I57").Value = 8665.333000000001
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 25995.999
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -25436.999
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 825
$ws.Range("I129").Value = 825
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2475
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2525
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2351.0876
$ws.Range("I131").Value = 1541.25
$ws.Range("J131").Value = 2412.2075
$ws.Range("K131").Value = 4623.75
$ws.Range("L131").Value = 7236.622499999999
$ws.Range("M131").Value = 416.25
$ws.Range("N131").Value = -17316.6225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 25000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 25000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H67").Value = 25000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 25000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1598.4615
$ws.Range("I107").Value = 1648.3334
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1648.3334
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 271.6666
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 41999.69
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 41999.69
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 41999.69
$ws.Range("N109").Value = -44079.69

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 10565.333
$ws.Range("I16").Value = 9898.5
$ws.Range("J16").Value = 11899
$ws.Range("K16").Value = 9898.5
$ws.Range("L16").Value = 11899
$ws.Range("M16").Value = -9728.5
$ws.Range("N16").Value = -12239

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1184.8572
$ws.Range("I22").Value = 799.5
$ws.Range("J22").Value = 1339
$ws.Range("K22").Value = 799.5
$ws.Range("L22").Value = 1339
$ws.Range("M22").Value = -504.5
$ws.Range("N22").Value = -1929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1184.8572
$ws.Range("I27").Value = 799.5
$ws.Range("J27").Value = 1339
$ws.Range("K27").Value = 799.5
$ws.Range("L27").Value = 1339
$ws.Range("M27").Value = -692.5
$ws.Range("N27").Value = -1553

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2224.875
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 2419.8
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 2419.8
$ws.Range("M46").Value = -1712
$ws.Range("N46").Value = -2795.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1049.5
$ws.Range("I55").Value = 1109.6
$ws.Range("J55").Value = 749
$ws.Range("K55").Value = 1109.6
$ws.Range("L55").Value = 749
$ws.Range("M55").Value = -936.5999999999999
$ws.Range("N55").Value = -1095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4950
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 4950
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 4950
$ws.Range("N61").Value = -5354

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4950
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4950
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4950
$ws.Range("N113").Value = -9290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 8000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 24000
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -21450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 5250
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 10000
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 10000
$ws.Range("M55").Value = -223
$ws.Range("N55").Value = -10554

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1459.8
$ws.Range("I122").Value = 1275.9412
$ws.Range("J122").Value = 2501.6667
$ws.Range("K122").Value = 3827.8236
$ws.Range("L122").Value = 7505.000100000001
$ws.Range("M122").Value = -1377.8236
$ws.Range("N122").Value = -12405.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5454.1
$ws.Range("I136").Value = 3950.4443
$ws.Range("J136").Value = 18987
$ws.Range("K136").Value = 11851.3329
$ws.Range("L136").Value = 56961
$ws.Range("M136").Value = -9301.332900000001
$ws.Range("N136").Value = -62061
